$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 127.0108868979071
$ws.Range("C2").Value = 6136
$ws.Range("D2").Value = 0.020699297082449
$ws.Range("E2").Value = 127.0108868979071
$ws.Range("B3").Value = -61.24612018567306
$ws.Range("C3").Value = 732
$ws.Range("D3").Value = -0.08366956309518178
$ws.Range("E3").Value = -61.24612018567306
$ws.Range("B4").Value = -212.3418476210031
$ws.Range("C4").Value = 3904
$ws.Range("D4").Value = -0.05439084211603563
$ws.Range("E4").Value = -212.3418476210031
$ws.Range("B5").Value = -83.25904031629506
$ws.Range("C5").Value = 1444
$ws.Range("D5").Value = -0.05765861517748965
$ws.Range("E5").Value = -83.25904031629506
$ws.Range("B6").Value = -381.8469510135291
$ws.Range("C6").Value = 3744
$ws.Range("D6").Value = -0.1019890360613059
$ws.Range("E6").Value = -381.8469510135291
$ws.Range("B7").Value = -226.9751793531174
$ws.Range("C7").Value = 3392
$ws.Range("D7").Value = -0.06691485240363132
$ws.Range("E7").Value = -226.9751793531174
$ws.Range("B8").Value = -156.8876156948677
$ws.Range("C8").Value = 3296
$ws.Range("D8").Value = -0.0475993979656759
$ws.Range("E8").Value = -156.8876156948677
$ws.Range("B9").Value = -82.32936737153875
$ws.Range("C9").Value = 3924
$ws.Range("D9").Value = -0.02098098047184983
$ws.Range("E9").Value = -82.32936737153875
$ws.Range("B10").Value = 37.18195107458128
$ws.Range("C10").Value = 1632
$ws.Range("D10").Value = 0.02278305825648362
$ws.Range("E10").Value = 37.18195107458128
$ws.Range("B11").Value = -149.662743086256
$ws.Range("C11").Value = 2608
$ws.Range("D11").Value = -0.05738602112203069
$ws.Range("E11").Value = -149.662743086256
$ws.Range("B12").Value = -260.2462975575222
$ws.Range("C12").Value = 1844
$ws.Range("D12").Value = -0.1411313978077669
$ws.Range("E12").Value = -260.2462975575222
$ws.Range("B13").Value = -125.4879013688354
$ws.Range("C13").Value = 888
$ws.Range("D13").Value = -0.141315204244184
$ws.Range("E13").Value = -125.4879013688354
$ws.Range("B14").Value = 45.489393271398
$ws.Range("C14").Value = 1464
$ws.Range("D14").Value = 0.03107198993947951
$ws.Range("E14").Value = 45.489393271398
$ws.Range("B15").Value = -720.0590425069661
$ws.Range("C15").Value = 7348
$ws.Range("D15").Value = -0.09799388166942924
$ws.Range("E15").Value = -720.0590425069661
$ws.Range("B16").Value = 47.70895354011265
$ws.Range("C16").Value = 2560
$ws.Range("D16").Value = 0.0186363099766065
$ws.Range("E16").Value = 47.70895354011265
$ws.Range("B17").Value = -594.9757865125223
$ws.Range("C17").Value = 3896
$ws.Range("D17").Value = -0.152714524258861
$ws.Range("E17").Value = -594.9757865125223
$ws.Range("B18").Value = -32.93947489898112
$ws.Range("C18").Value = 976
$ws.Range("D18").Value = -0.03374946198666098
$ws.Range("E18").Value = -32.93947489898112
$ws.Range("B19").Value = -450.8703413625891
$ws.Range("C19").Value = 4636
$ws.Range("D19").Value = -0.09725417199365598
$ws.Range("E19").Value = -450.8703413625891
$ws.Range("B20").Value = -112.7130100074113
$ws.Range("C20").Value = 1188
$ws.Range("D20").Value = -0.09487627105000955
$ws.Range("E20").Value = -112.7130100074113
$ws.Range("B21").Value = -50.26298722493762
$ws.Range("C21").Value = 2252
$ws.Range("D21").Value = -0.0223192660856739
$ws.Range("E21").Value = -50.26298722493762
$ws.Range("B22").Value = -109.5953780685138
$ws.Range("C22").Value = 2624
$ws.Range("D22").Value = -0.04176653127611044
$ws.Range("E22").Value = -109.5953780685138
$ws.Range("B23").Value = 19.07818271353684
$ws.Range("C23").Value = 1228
$ws.Range("D23").Value = 0.01553597940841762
$ws.Range("E23").Value = 19.07818271353684
$ws.Range("B24").Value = -141.1988355505335
$ws.Range("C24").Value = 1956
$ws.Range("D24").Value = -0.07218754373749157
$ws.Range("E24").Value = -141.1988355505335
$ws.Range("B25").Value = 20.9436434891197
$ws.Range("C25").Value = 1408
$ws.Range("D25").Value = 0.01487474679624979
$ws.Range("E25").Value = 20.9436434891197
$ws.Range("B26").Value = -282.5282828522926
$ws.Range("C26").Value = 2940
$ws.Range("D26").Value = -0.09609805539193625
$ws.Range("E26").Value = -282.5282828522926
$ws.Range("B27").Value = -59.25921405929971
$ws.Range("C27").Value = 1244
$ws.Range("D27").Value = -0.0476360241634242
$ws.Range("E27").Value = -59.25921405929971
$ws.Range("B28").Value = -226.8913927195349
$ws.Range("C28").Value = 1628
$ws.Range("D28").Value = -0.1393681773461516
$ws.Range("E28").Value = -226.8913927195349
$ws.Range("B29").Value = -448.4560768004027
$ws.Range("C29").Value = 2496
$ws.Range("D29").Value = -0.1796699025642639
$ws.Range("E29").Value = -448.4560768004027
$ws.Range("B30").Value = 191.5671265455729
$ws.Range("C30").Value = 6580
$ws.Range("D30").Value = 0.02911354506771625
$ws.Range("E30").Value = 191.5671265455729
$ws.Range("B31").Value = -85.95745261293301
$ws.Range("C31").Value = 784
$ws.Range("D31").Value = -0.1096396079246594
$ws.Range("E31").Value = -85.95745261293301
$ws.Range("B32").Value = -102.6527857225576
$ws.Range("C32").Value = 3068
$ws.Range("D32").Value = -0.03345918700213742
$ws.Range("E32").Value = -102.6527857225576
$ws.Range("B33").Value = -128.8681815330119
$ws.Range("C33").Value = 1716
$ws.Range("D33").Value = -0.0750980078863706
$ws.Range("E33").Value = -128.8681815330119
$ws.Range("B34").Value = -363.6641642799185
$ws.Range("C34").Value = 2900
$ws.Range("D34").Value = -0.1254014359585926
$ws.Range("E34").Value = -363.6641642799185
$ws.Range("B35").Value = -124.3523736193784
$ws.Range("C35").Value = 1980
$ws.Range("D35").Value = -0.06280422910069618
$ws.Range("E35").Value = -124.3523736193784
$ws.Range("B36").Value = 58.44100642036419
$ws.Range("C36").Value = 1692
$ws.Range("D36").Value = 0.03453960190328853
$ws.Range("E36").Value = 58.44100642036419
$ws.Range("B37").Value = 48.02254751488483
$ws.Range("C37").Value = 2892
$ws.Range("D37").Value = 0.01660530688619808
$ws.Range("E37").Value = 48.02254751488483
$ws.Range("B38").Value = -88.97232123756133
$ws.Range("C38").Value = 736
$ws.Range("D38").Value = -0.1208863060292953
$ws.Range("E38").Value = -88.97232123756133
$ws.Range("B39").Value = -504.1183246364689
$ws.Range("C39").Value = 3320
$ws.Range("D39").Value = -0.1518428688664063
$ws.Range("E39").Value = -504.1183246364689
$ws.Range("B40").Value = -362.1260193438251
$ws.Range("C40").Value = 3424
$ws.Range("D40").Value = -0.1057611037803227
$ws.Range("E40").Value = -362.1260193438251
$ws.Range("B41").Value = -65.5626212472756
$ws.Range("C41").Value = 2140
$ws.Range("D41").Value = -0.03063673890059607
$ws.Range("E41").Value = -65.5626212472756
